$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.0.0"
$meta.Range("B8").Value = "2025-06-05T14:31:57+02:00"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI4").Value = ""
$elements.Range("AI6").Value = ""
